# Add four new legend-player rows (88-91) with their stats, update the
# sheet's view/selection to reflect the new bottom of the data, and set
# the page setup (paper size / orientation) for printing.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New rows: 88 (Charles Barkley) .. 91 (Elton Brand) ------------------
$newPlayers = @(
    @{ Row = 88; Name = "Charles Barkley"; B = 1073; C = 0.541; D = 0.266; E = 0.684; F = 12546; G = 42158; H = 23757 },
    @{ Row = 89; Name = "Nick Anderson";   B = 800;  C = 0.446; D = 0.356; E = 0.667; F = 4064;  G = 2087;  H = 11529 },
    @{ Row = 90; Name = "Grant Long";      B = 1003; C = 0.467; D = 0.283; E = 0.761; F = 6154;  G = 1716;  H = 9518  },
    @{ Row = 91; Name = "Elton Brand";     B = 1058; C = 0.5;   D = 0.095; E = 0.736; F = 9040;  G = 2184;  H = 16827 }
)

foreach ($p in $newPlayers) {
    $r = $p.Row
    $ws.Range("A$r").Value = $p.Name
    $ws.Range("B$r").Value = $p.B
    $ws.Range("C$r").Value = $p.C
    $ws.Range("D$r").Value = $p.D
    $ws.Range("E$r").Value = $p.E
    $ws.Range("F$r").Value = $p.F
    $ws.Range("G$r").Value = $p.G
    $ws.Range("H$r").Value = $p.H
}

# --- View state: scroll the window down and move the selection ----------
$excel.ActiveWindow.ScrollRow = 85
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("K96").Select()

# --- Page setup: A4, portrait -------------------------------------------
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

Write-Output "added legend rows 88-91, updated view and page setup"
